$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2  = 4
  3  = -1
  4  = -1
  6  = 2
  7  = 2
  8  = -3
  9  = 2
  10 = -5
  11 = -3
  12 = 2
  13 = 4
  14 = -5
  15 = 1
  17 = -4
  19 = -2
  21 = 2
  22 = -1
  23 = -3
  25 = 3
  26 = 2
  27 = -2
  28 = -1
  29 = 2
  31 = -1
}

foreach ($row in $values.Keys) {
  $ws.Range("F$row").Value = $values[$row]
}
